$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$g2 = "<get>`n  <filter>`n    <components xmlns=""http://openconfig.net/yang/platform"">`n      <component>`n        <name>Waveserver-Ai</name>`n        <state>`n          <part-no></part-no>`n        </state>`n      </component>`n    </components>`n  </filter>`n</get>"

$j2 = "<?xml version=""1.0"" encoding=""UTF-8""?>`n<rpc-reply message-id=""urn:uuid:e6891d16-81b1-44eb-8d42-93065de7beed""`n xmlns:ncx=""http://netconfcentral.org/ns/yuma-ncx""`n ncx:last-modified=""2020-10-07T13:51:28Z"" ncx:etag=""814""`n xmlns=""urn:ietf:params:xml:ns:netconf:base:1.0"">`n <data>`n  <components xmlns=""http://openconfig.net/yang/platform"">`n   <component>`n    <name>Waveserver-Ai</name>`n    <state>`n     <part-no>186-1010-900</part-no>`n    </state>`n   </component>`n  </components>`n </data>`n</rpc-reply>"

$ws.Range("G2").Value = $g2
$ws.Range("J2").Value = $j2

$wb.Save()
